# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates the DAMSLTag (column I) and DialogAct (column J) values for a set of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 3;   DamslTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 7;   DamslTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 15;  DamslTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 34;  DamslTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 35;  DamslTag = "sd"; DialogAct = "Statement-non-opinion" }
    @{ Row = 53;  DamslTag = "%"; DialogAct = "Uninterpretable" }
    @{ Row = 72;  DamslTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 83;  DamslTag = "aa"; DialogAct = "Agree/Accept" }
    @{ Row = 92;  DamslTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 103; DamslTag = "sv"; DialogAct = "Statement-opinion" }
    @{ Row = 104; DamslTag = "%"; DialogAct = "Uninterpretable" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.DamslTag
    $ws.Cells.Item($u.Row, 10).Value = $u.DialogAct
}
